$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 233

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 142

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 105

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 76

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 74

$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
